$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ToDO")

# Row 7: "Werkelijke tijd" typo fix "3uur" -> "4uur"
$ws.Range("C7").Value = "4uur"

# Row 11: new task "helpen trigger" (done Thu 28/03)
$ws.Range("A11").Value = "helpen trigger"
$ws.Range("B11").Value = "20 min"
$ws.Range("C11").Value = "3 uur"

# Copy date formatting from the cell above so the new date cell matches
# the existing "Voltooid" column formatting instead of minting a new style.
$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").Value = 41361

$ws.Range("E11").Value = 1
$ws.Range("F11").Value = "Wouter"
$ws.Range("G11").Value = "In Process"
$ws.Range("H11").Value = "sql"
